# Add 2022-Q4 data:
#   1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3" sheet,
#      seeded as a structural copy of "2022-Q3" (same headers/styles/text types),
#      then overwrite its fund-level numbers with the new quarter's figures.
#   2. Update the "总计" (Total) summary sheet: push the existing quarterly
#      rows down by one and insert the new 2022-Q4 summary row at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q4" sheet
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3) | Out-Null          # new copy is placed immediately before $q3
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Fund-level figures for 2022-Q4 (columns D/E/F/G are stored as text, like the
# rest of the workbook, so force text with a leading apostrophe).
$q4.Cells.Item(2, 4).Value = "'2.81"
$q4.Cells.Item(2, 5).Value = "'93.63"
$q4.Cells.Item(2, 6).Value = "'2.87"
$q4.Cells.Item(2, 7).Value = "'0.0806"

$q4.Cells.Item(3, 4).Value = "'0.36"
$q4.Cells.Item(3, 5).Value = "'93.63"
$q4.Cells.Item(3, 6).Value = "'2.87"
$q4.Cells.Item(3, 7).Value = "'0.0103"

# ---------------------------------------------------------------------------
# 2. "总计" summary sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Duplicate the last row's formatting into the newly-needed row 10 so the
# appended 2020-Q4 row picks up the same styling as every other data row.
$total.Range("A9:D9").Copy($total.Range("A10:D10")) | Out-Null

$rows = @(
    @(2,  0, "2022-Q4", 2, 0.09),
    @(3,  1, "2022-Q3", 2, 0.09),
    @(4,  2, "2022-Q2", 3, 0.22),
    @(5,  3, "2022-Q1", 2, 0.13),
    @(6,  4, "2021-Q4", 1, 0.11),
    @(7,  5, "2021-Q3", 2, 0.14),
    @(8,  6, "2021-Q2", 1, 0.14),
    @(9,  7, "2021-Q1", 1, 0.19),
    @(10, 8, "2020-Q4", 1, 0.31)
)

foreach ($row in $rows) {
    $r = $row[0]
    $total.Cells.Item($r, 1).Value = $row[1]
    $total.Cells.Item($r, 2).Value = $row[2]
    $total.Cells.Item($r, 3).Value = $row[3]
    $total.Cells.Item($r, 4).Value = $row[4]
}
